$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The B column held dates typed as "14/4/20NN" text strings (day/month/year,
# one per row). Replace them with genuine Excel date serials - 14 Apr 2001
# through 24 Apr 2001 - stored as numbers and displayed with a date format,
# instead of as text.
$dateSerials = 36995, 36996, 36997, 36998, 36999, 37000, 37001, 37002, 37003, 37004, 37005
for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $dateSerials[$i]
}
$ws.Range("B1:B11").NumberFormat = "m/d/yy;@"

# Set the print orientation to portrait.
$ws.PageSetup.Orientation = 1

# Leave the on-screen selection on the single cell B8.
[void]$ws.Range("B8").Select()
